# Correct the CCS subsidy duration used in the 45Q credit-value calculation
# from 10 years (erroneous) to 12 years, matching the IRA credit duration.

$wb = $excel.ActiveWorkbook

$wsAbout = $wb.Worksheets.Item("About")
$wsElec  = $wb.Worksheets.Item("Electricity Calculations")

# --- Core data fix -------------------------------------------------------
# "45Q Credit Duration" (B3) was 10, should be 12 years.
# B4 ("Credit Value" = B1*B3/B2) and everything downstream on the BCS sheet
# that references it recalculates automatically.
$wsElec.Range("B3").Value = 12

# --- Formatting cleanup ----------------------------------------------------
# Clear the (no-op) font style that had been applied to the note in About!A7.
$wsAbout.Range("A7").ClearFormats()

# --- Restore the UI selection state left after the edit --------------------
$wsElec.Activate()
$wsElec.Range("B4").Select()

$wsAbout.Activate()
$wsAbout.Range("B11").Select()
